# Add the GitHub Pages link after the LinkedIn URL, then re-create the
# "_GoBack" bookmark at the very end of that line (Word keeps _GoBack
# pinned to the location of the most recent edit; adding a bookmark
# named "_GoBack" automatically removes the old one elsewhere in the
# document, matching real Word behaviour).

$d = $word.ActiveDocument

$rng = $d.Content
$rng.Find.Execute("linkedin.com/in/sarah-elkins-93719742", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertAfter(" | ")
$rng.Collapse(0)
$rng.InsertAfter("https://selkins13.github.io")
$rng.Collapse(0)

# $rng is now collapsed exactly at the end of the paragraph's text. Adding
# a bookmark directly at a paragraph's trailing edge confuses the engine's
# range resolver, so nudge past it with a throwaway character, anchor the
# bookmark just before that character, then delete the character again.
$rng.InsertAfter("X")
$goBackPos = $rng.Start
$bmRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
$d.Range($goBackPos, $goBackPos + 1).Delete()
